# Add the 2023 column (column S) to the "hotels and restaurants" sheet,
# mirroring the existing 2022 column's (column R) formatting for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the formatting of the 2022 column (R3:R14) across to the new
# 2023 column (S3:S14) before writing the new figures.
$ws.Range("R3:R14").Copy()
$ws.Range("S3:S14").PasteSpecial(-4122)

# Header
$ws.Range("S3").Value = 2023

# Data rows (same order as the existing indicator rows 4-14)
$ws.Range("S4").Value = 33.9
$ws.Range("S5").Value = 33.9
$ws.Range("S6").Value = 854
$ws.Range("S7").Value = 842
$ws.Range("S8").Value = 649.16999999999996
$ws.Range("S9").Value = 24.2
$ws.Range("S10").Value = 6.6
$ws.Range("S11").Value = 9.6999999999999993
$ws.Range("S12").Value = 0.8
$ws.Range("S13").Value = 24.1
$ws.Range("S14").Value = "_"

# Match the author's final selection state
$ws.Range("E22").Select()
